$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano test statistics (DM_Stat, column C) and
# p-values (P_Value, column D) for rows 2-11.

$ws.Range("C2").Value = 0.6827182936900972
$ws.Range("D2").Value = 0.5019112995627282

$ws.Range("C3").Value = -0.6034871086408529
$ws.Range("D3").Value = 0.5523566792176868

$ws.Range("C4").Value = -0.2095395701032416
$ws.Range("D4").Value = 0.8359552927625062

$ws.Range("C5").Value = 0.49460756977287
$ws.Range("D5").Value = 0.6257803174281276

$ws.Range("C6").Value = -1.330260272875008
$ws.Range("D6").Value = 0.1970588732705119

$ws.Range("C7").Value = -0.69909142913724
$ws.Range("D7").Value = 0.4918204365806202

$ws.Range("C8").Value = -0.002822321873415189
$ws.Range("D8").Value = 0.9977735521300291

$ws.Range("C9").Value = 0.3814525067794956
$ws.Range("D9").Value = 0.7065267412498839

$ws.Range("C10").Value = 1.972141351171651
$ws.Range("D10").Value = 0.06129442439186739

$ws.Range("C11").Value = 0.5211460389418854
$ws.Range("D11").Value = 0.6074716359589756
